$wb = $excel.ActiveWorkbook

# ---- Year sheets 2006-2010: clear Alaska (FIPS 2, row 3) and Hawaii (FIPS 15, row 13) ----
# and remove trailing placeholder rows with no state name (FIPS 66/72/78).
$yearSheets = @("2006", "2007", "2008", "2009", "2010")
foreach ($name in $yearSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C3:E3").ClearContents()
    $ws.Range("C13:E13").ClearContents()
}

# 2006 sheet only has 2 trailing placeholder rows (66 is absent); others have 3.
$ws2006 = $wb.Worksheets.Item("2006")
$ws2006.Rows.Item(54).Delete()
$ws2006.Rows.Item(53).Delete()

$otherYearSheets = @("2007", "2008", "2009", "2010")
foreach ($name in $otherYearSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(55).Delete()
    $ws.Rows.Item(54).Delete()
    $ws.Rows.Item(53).Delete()
}

# ---- Aggregate sheet: recompute combined counts now that Alaska/Hawaii are empty ----
$agg = $wb.Worksheets.Item("Aggregate")

$agg.Cells.Item(2, 1).Value = 4.0
$agg.Cells.Item(2, 2).Value = "Arizona"
$agg.Cells.Item(2, 3).Value = 42621.918836024626
$agg.Cells.Item(2, 4).Value = 2802421.7470672023
$agg.Cells.Item(2, 5).Value = 15.208959493919652

$agg.Cells.Item(3, 1).Value = 6.0
$agg.Cells.Item(3, 2).Value = "California"
$agg.Cells.Item(3, 3).Value = 156599.0264181489
$agg.Cells.Item(3, 4).Value = 16850453.394726697
$agg.Cells.Item(3, 5).Value = 9.29346070101331

$agg.Cells.Item(4, 1).Value = 9.0
$agg.Cells.Item(4, 2).Value = "Connecticut"
$agg.Cells.Item(4, 3).Value = 32938.5229378051
$agg.Cells.Item(4, 4).Value = 2734477.7688722718
$agg.Cells.Item(4, 5).Value = 12.045635664973537

$agg.Cells.Item(5, 1).Value = 11.0
$agg.Cells.Item(5, 2).Value = "District of Columbia"
$agg.Cells.Item(5, 3).Value = 3184.446344618885
$agg.Cells.Item(5, 4).Value = 179492.85599165817
$agg.Cells.Item(5, 5).Value = 17.74135425627681

$agg.Cells.Item(6, 1).Value = 13.0
$agg.Cells.Item(6, 2).Value = "Geogia"
$agg.Cells.Item(6, 3).Value = 94786.03809289802
$agg.Cells.Item(6, 4).Value = 10458073.543508092
$agg.Cells.Item(6, 5).Value = 9.063431969431596

$agg.Cells.Item(7, 1).Value = 17.0
$agg.Cells.Item(7, 2).Value = "Illinois"
$agg.Cells.Item(7, 3).Value = 37799.24289699028
$agg.Cells.Item(7, 4).Value = 5673571.265652148
$agg.Cells.Item(7, 5).Value = 6.662336846956915

$agg.Cells.Item(8, 1).Value = 18.0
$agg.Cells.Item(8, 2).Value = "Indiana"
$agg.Cells.Item(8, 3).Value = 105218.93189138013
$agg.Cells.Item(8, 4).Value = 6936762.371679282
$agg.Cells.Item(8, 5).Value = 15.168305652354107

$agg.Cells.Item(9, 1).Value = 19.0
$agg.Cells.Item(9, 2).Value = "Iowa"
$agg.Cells.Item(9, 3).Value = 11510.486502884885
$agg.Cells.Item(9, 4).Value = 1829733.9661702944
$agg.Cells.Item(9, 5).Value = 6.290797851327419

$agg.Cells.Item(10, 1).Value = 20.0
$agg.Cells.Item(10, 2).Value = "Kansas"
$agg.Cells.Item(10, 3).Value = 27509.487892604753
$agg.Cells.Item(10, 4).Value = 3059760.3617240055
$agg.Cells.Item(10, 5).Value = 8.990732815789757

$agg.Cells.Item(11, 1).Value = 22.0
$agg.Cells.Item(11, 2).Value = "Louisiana"
$agg.Cells.Item(11, 3).Value = 5378.786540781787
$agg.Cells.Item(11, 4).Value = 931966.3594570106
$agg.Cells.Item(11, 5).Value = 5.771438514063552

$agg.Cells.Item(12, 1).Value = 23.0
$agg.Cells.Item(12, 2).Value = "Maine"
$agg.Cells.Item(12, 3).Value = 6662.149701015229
$agg.Cells.Item(12, 4).Value = 722763.1937653258
$agg.Cells.Item(12, 5).Value = 9.217610634415294

$agg.Cells.Item(13, 1).Value = 24.0
$agg.Cells.Item(13, 2).Value = "Maryland"
$agg.Cells.Item(13, 3).Value = 64870.5736086545
$agg.Cells.Item(13, 4).Value = 5816583.904105316
$agg.Cells.Item(13, 5).Value = 11.152692831073782

$agg.Cells.Item(14, 1).Value = 26.0
$agg.Cells.Item(14, 2).Value = "Michigan"
$agg.Cells.Item(14, 3).Value = 126101.97053771566
$agg.Cells.Item(14, 4).Value = 10491065.363787048
$agg.Cells.Item(14, 5).Value = 12.0199394594369

$agg.Cells.Item(15, 1).Value = 28.0
$agg.Cells.Item(15, 2).Value = "Mississippi"
$agg.Cells.Item(15, 3).Value = 18264.008952066375
$agg.Cells.Item(15, 4).Value = 1300916.7417331804
$agg.Cells.Item(15, 5).Value = 14.039337312035567

$agg.Cells.Item(16, 1).Value = 29.0
$agg.Cells.Item(16, 2).Value = "Missouri"
$agg.Cells.Item(16, 3).Value = 46410.034572168624
$agg.Cells.Item(16, 4).Value = 3600272.063188931
$agg.Cells.Item(16, 5).Value = 12.890702079625912

$agg.Cells.Item(17, 1).Value = 30.0
$agg.Cells.Item(17, 2).Value = "Montana"
$agg.Cells.Item(17, 3).Value = 3295.889765187216
$agg.Cells.Item(17, 4).Value = 768012.0066596719
$agg.Cells.Item(17, 5).Value = 4.291456040540417

$agg.Cells.Item(18, 1).Value = 31.0
$agg.Cells.Item(18, 2).Value = "Nebraska"
$agg.Cells.Item(18, 3).Value = 18262.226573581298
$agg.Cells.Item(18, 4).Value = 2014605.1364890952
$agg.Cells.Item(18, 5).Value = 9.064916118206348

$agg.Cells.Item(19, 1).Value = 33.0
$agg.Cells.Item(19, 2).Value = "New Hampshire"
$agg.Cells.Item(19, 3).Value = 9423.24961501254
$agg.Cells.Item(19, 4).Value = 788301.5634397555
$agg.Cells.Item(19, 5).Value = 11.953863917120971

$agg.Cells.Item(20, 1).Value = 34.0
$agg.Cells.Item(20, 2).Value = "New Jersey"
$agg.Cells.Item(20, 3).Value = 51471.692726499445
$agg.Cells.Item(20, 4).Value = 5274309.777003881
$agg.Cells.Item(20, 5).Value = 9.758943805484705

$agg.Cells.Item(21, 1).Value = 35.0
$agg.Cells.Item(21, 2).Value = "New Mexico"
$agg.Cells.Item(21, 3).Value = 8857.142946737338
$agg.Cells.Item(21, 4).Value = 1327496.3453727677
$agg.Cells.Item(21, 5).Value = 6.6720657858008705

$agg.Cells.Item(22, 1).Value = 36.0
$agg.Cells.Item(22, 2).Value = "New York"
$agg.Cells.Item(22, 3).Value = 221226.36562273267
$agg.Cells.Item(22, 4).Value = 15027480.545143578
$agg.Cells.Item(22, 5).Value = 14.721454135851552

$agg.Cells.Item(23, 1).Value = 39.0
$agg.Cells.Item(23, 2).Value = "Ohio"
$agg.Cells.Item(23, 3).Value = 71567.97481901593
$agg.Cells.Item(23, 4).Value = 4755244.714276284
$agg.Cells.Item(23, 5).Value = 15.050324246016872

$agg.Cells.Item(24, 1).Value = 40.0
$agg.Cells.Item(24, 2).Value = "Oklahoma"
$agg.Cells.Item(24, 3).Value = 24627.721388718226
$agg.Cells.Item(24, 4).Value = 2285659.3166476926
$agg.Cells.Item(24, 5).Value = 10.774887232467767

$agg.Cells.Item(25, 1).Value = 41.0
$agg.Cells.Item(25, 2).Value = "Oregon"
$agg.Cells.Item(25, 3).Value = 8328.024176806282
$agg.Cells.Item(25, 4).Value = 752767.6262806169
$agg.Cells.Item(25, 5).Value = 11.063207138641959

$agg.Cells.Item(26, 1).Value = 42.0
$agg.Cells.Item(26, 2).Value = "Pennsylvania"
$agg.Cells.Item(26, 3).Value = 62291.57800717453
$agg.Cells.Item(26, 4).Value = 4733924.715291146
$agg.Cells.Item(26, 5).Value = 13.158548509646815

$agg.Cells.Item(27, 1).Value = 44.0
$agg.Cells.Item(27, 2).Value = "Rhode Island"
$agg.Cells.Item(27, 3).Value = 5476.426426490039
$agg.Cells.Item(27, 4).Value = 384116.9213563548
$agg.Cells.Item(27, 5).Value = 14.257186085820528

$agg.Cells.Item(28, 1).Value = 48.0
$agg.Cells.Item(28, 2).Value = "Texas"
$agg.Cells.Item(28, 3).Value = 381999.04549088026
$agg.Cells.Item(28, 4).Value = 22992022.96837785
$agg.Cells.Item(28, 5).Value = 16.61441648767765

$agg.Cells.Item(29, 1).Value = 49.0
$agg.Cells.Item(29, 2).Value = "Utah"
$agg.Cells.Item(29, 3).Value = 30221.406332381404
$agg.Cells.Item(29, 4).Value = 2902955.358223645
$agg.Cells.Item(29, 5).Value = 10.410565304343592

$agg.Cells.Item(30, 1).Value = 50.0
$agg.Cells.Item(30, 2).Value = "Vermont"
$agg.Cells.Item(30, 3).Value = 6498.410650110485
$agg.Cells.Item(30, 4).Value = 563279.5559912089
$agg.Cells.Item(30, 5).Value = 11.536741536225586

$agg.Cells.Item(31, 1).Value = 53.0
$agg.Cells.Item(31, 2).Value = "Washington"
$agg.Cells.Item(31, 3).Value = 18647.25128699706
$agg.Cells.Item(31, 4).Value = 2752372.813938709
$agg.Cells.Item(31, 5).Value = 6.774972922477175

$agg.Cells.Item(32, 1).Value = 54.0
$agg.Cells.Item(32, 2).Value = "West Virginia"
$agg.Cells.Item(32, 3).Value = 3846.644645155366
$agg.Cells.Item(32, 4).Value = 325031.2108367552
$agg.Cells.Item(32, 5).Value = 11.834693152244132

$agg.Cells.Item(33, 1).Value = 55.0
$agg.Cells.Item(33, 2).Value = "Wisconsin"
$agg.Cells.Item(33, 3).Value = 14403.531708608236
$agg.Cells.Item(33, 4).Value = 1174446.719861135
$agg.Cells.Item(33, 5).Value = 12.264099737373606

# Remove the two now-obsolete trailing rows (sheet shrinks from 35 to 33 total rows).
$agg.Rows.Item(35).Delete()
$agg.Rows.Item(34).Delete()

